$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C27").Value = "Java 8 continued"
Write-Host "done"
